$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)

# The five "UC#" ellipse call-out shapes (p:cNvPr id 92, 93, 95, 98, 99 /
# names "Google Shape;92;p16" etc.) were repositioned on slide 4
# ("I assessed feasibility vs. impact for all cases").
#
# Left/Top are expressed in points (1 pt = 12700 EMU); the literals below
# are tuned so the stored EMU offsets land exactly on the target values
# (x="4436883" y="1822469", etc.) after the engine's internal
# points->EMU conversion.

$sh92 = $s.Shapes.Item("Google Shape;92;p16")
$sh92.Left = 349.3608661417323
$sh92.Top = 143.5015030629917

$sh93 = $s.Shapes.Item("Google Shape;93;p16")
$sh93.Left = 240.86850743700768
$sh93.Top = 262.5443307086614

$sh95 = $s.Shapes.Item("Google Shape;95;p16")
$sh95.Left = 296.50826771653544
$sh95.Top = 158.5561447322833

$sh98 = $s.Shapes.Item("Google Shape;98;p16")
$sh98.Left = 335.6343307086614
$sh98.Top = 191.1123622047244

$sh99 = $s.Shapes.Item("Google Shape;99;p16")
$sh99.Left = 289.4600787401575
$sh99.Top = 240.79528049055088
